# Update the "matrices" scores (column F) and re-sort the whole table by
# that column (descending), refreshing prolificid/name/gender/index to
# follow their respective person as the ranking shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13 (columns A..H):
#   A = row index (0-based, unchanged)
#   B = index
#   C = prolificid
#   D = name
#   E = gender
#   F = matrices (the values that actually changed)
#   G = race
#   H = mat_rank (1-based rank, unchanged)
$rows = @(
    @(0, 0,  "5c27de12a2b00a00018b2c16", "Ankai",    "male",   14.36475064273752, "Asian", 1),
    @(1, 2,  "5f2c1a97a6809c060fec8820", "Maggie",   "female", 13.4427811560038,  "Asian", 2),
    @(2, 1,  "60b1742bce2b39e0f1d19a1a", "Sabrina",  "female", 13.32257368402617, "Asian", 3),
    @(3, 3,  "60bd88b8fc436774352f53b9", "Annes",    "female", 13.02548504840682, "Asian", 4),
    @(4, 4,  "5f7cbf8a2fe61814cae2ce8b", "Aalap",    "male",   12.16366162123603, "Asian", 5),
    @(5, 6,  "60b7cd4be2d4cc6bb252e016", "Chris",    "male",   10.35758251781631, "Asian", 6),
    @(6, 7,  "5ff3974450a7199965624df7", "Anh",      "male",   10.23661900101856, "Asian", 7),
    @(7, 10, "60a71d27a66fac796ad4de6f", "Jennifer", "female", 8.201924197465678, "Asian", 8),
    @(8, 13, "5697d4ae7183b8000d0fc201", "Tu",       "male",   5.441561929436489, "Asian", 9),
    @(9, 14, "60186dc2cc1aa8103499603a", "Emily",    "female", 2.330660576781288, "Asian", 10),
    @(10,15, "60b76ee2219ac1ce25ccea43", "Richie",   "male",   2.005372734962068, "Asian", 11),
    @(11,16, "60863a15760523386e761cfb", "Roshni",   "female", 1.34066941120993,  "Asian", 12)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
